$d = $word.ActiveDocument

# 1) "Oracle ADF,IBM WebSphere,ASP.NET WebForms,Custom banking front ends"
#    -> "Oracle, java,asp.net webforms,Custom banking front ends"
$d.Content.Find.Execute(
    "Oracle ADF,IBM WebSphere,ASP.NET WebForms,Custom banking front ends",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Oracle, java,asp.net webforms,Custom banking front ends", 2) | Out-Null

# 2) "Harness is not responsible for:" -> "and Harness is not responsible for:"
$d.Content.Find.Execute(
    "Harness is not responsible for:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and Harness is not responsible for:", 2) | Out-Null

# 3) After the "postback-driven frameworks" paragraph, insert:
#      - a blank paragraph
#      - a "Conclusion:" paragraph
#      - a paragraph with the closing findings text
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "postback-driven frameworks") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $insertPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)
    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $blankParaXml = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"default`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr></w:p>"
    $conclusionParaXml = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"default`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"default`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>Conclusion:</w:t></w:r></w:p>"
    $findingsParaXml = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"default`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"default`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>Based on the above findings, the current setup of Harness Cloud is not aligned with the testing requirements of banking applications. Most issues stem from domain-specific security restrictions and dynamic UI behavior, which cannot be addressed solely through Harness.</w:t></w:r></w:p>"
    $insertPoint.InsertXML($blankParaXml + $conclusionParaXml + $findingsParaXml)
}

# 4) Add explicit rFonts to the "Normal (Web)" style's run properties
#    (matching the document's rPrDefault: Times New Roman / SimSun)
$webStyle = $d.Styles.Item("Normal (Web)")
$webStyle.Font.NameAscii = "Times New Roman"
$webStyle.Font.NameOther = "Times New Roman"
$webStyle.Font.NameFarEast = "SimSun"
$webStyle.Font.NameBi = "Times New Roman"
